# Refresh the crypto price/volume snapshot (cols D = Price, E = Volume(1h))
# with the latest scraped values, per the GitHub Actions update commit.
#
# NumberFormat "@" is applied before writing column D prices because several
# of them are plain decimal-looking strings (e.g. "7.70", "27.78"); without
# forcing Text format first, Excel's COM layer auto-coerces such strings to
# numeric cells and silently drops trailing zeros (e.g. "7.70" -> 7.7),
# which would not match the source data (stored as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.553.93'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.674.42'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.44'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("E6").Value = '  +3.93%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.545'
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.672.32'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.359'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.78'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.160.80'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.403.46'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.654.02'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.70'
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '363.53'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("E22").Value = '  -3.30%  '
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  -4.22%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("E26").Value = '  -4.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.01'
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '556.36'
$ws.Range("E31").Value = '  -4.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.01'
$ws.Range("E32").Value = '  -2.17%  '
$ws.Range("E33").Value = '  -3.58%  '
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.52'
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("E39").Value = '  -3.91%  '
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.31'
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E42").Value = '  -4.07%  '
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E45").Value = '  -5.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.25'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0298'
$ws.Range("E47").Value = '  -5.88%  '
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '153.31'
$ws.Range("E49").Value = '  -2.66%  '
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("E51").Value = '  -2.77%  '
